# Updates the "Estado de Cuenta" worksheet:
#  - refresh the "Valor Mora" total and worker count
#  - rebuild the detail table with the new data set (adds a new worker,
#    a new "2508" period, and reorders rows by period)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header summary cells -----------------------------------------------
$ws.Range("E11").Value = 839520
$ws.Range("C13").Value = 4

# --- Detail table (B16:G29) ---------------------------------------------
# Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora, F=Valor Mora, G=Salario Basico
$data = @(
    @("CC", "1143352761", "ANDREA PAOLA VILLALOBOS SIMANCAS", "2503", 64000, 1600000),
    @("CC", "9294312",    "JHON JAIRO PAJARO ROJANO",         "2503", 56940, 1423500),
    @("CC", "1143352761", "ANDREA PAOLA VILLALOBOS SIMANCAS", "2504", 64000, 1600000),
    @("CC", "9294312",    "JHON JAIRO PAJARO ROJANO",         "2504", 56940, 1423500),
    @("CC", "73165996",   "WILSON ENRIQUE LORDUY LLERENA",    "2505", 56940, 1423500),
    @("CC", "1143352761", "ANDREA PAOLA VILLALOBOS SIMANCAS", "2505", 64000, 1600000),
    @("CC", "9294312",    "JHON JAIRO PAJARO ROJANO",         "2505", 56940, 1423500),
    @("CC", "1143352761", "ANDREA PAOLA VILLALOBOS SIMANCAS", "2506", 64000, 1600000),
    @("CC", "9294312",    "JHON JAIRO PAJARO ROJANO",         "2506", 56940, 1423500),
    @("CC", "1143352761", "ANDREA PAOLA VILLALOBOS SIMANCAS", "2507", 64000, 1600000),
    @("CC", "9294312",    "JHON JAIRO PAJARO ROJANO",         "2507", 56940, 1423500),
    @("CC", "1050969857", "JAINER DEIVISON CALDERON CARDENAS","2508", 56940, 1423500),
    @("CC", "1143352761", "ANDREA PAOLA VILLALOBOS SIMANCAS", "2508", 64000, 1600000),
    @("CC", "9294312",    "JHON JAIRO PAJARO ROJANO",         "2508", 56940, 1423500)
)

$startRow = 16
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]

    $ws.Cells.Item($r, 2).Value = [string]$rowVals[0]
    $ws.Cells.Item($r, 3).Value = [string]$rowVals[1]
    $ws.Cells.Item($r, 4).Value = [string]$rowVals[2]
    $ws.Cells.Item($r, 5).Value = [string]$rowVals[3]
    $ws.Cells.Item($r, 6).Value = $rowVals[4]
    $ws.Cells.Item($r, 7).Value = $rowVals[5]
}

# Column D best-fits to the new (longer) name so its width stays correct
$ws.Range("D16:D29").Columns.AutoFit() | Out-Null

$wb.Save()
